$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values with repulled data
$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -2
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 1
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = -9
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -5
$ws.Range("F27").Value = -3
$ws.Range("F28").Value = -1
$ws.Range("F29").Value = 7
$ws.Range("F30").Value = 3
$ws.Range("F35").Value = -4
$ws.Range("F37").Value = -5
$ws.Range("F40").Value = -9
$ws.Range("F42").Value = -1
$ws.Range("F43").Value = -2
$ws.Range("F45").Value = -6
$ws.Range("F46").Value = -1
$ws.Range("F47").Value = -5
$ws.Range("F49").Value = 1
